$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 103 -----------------------------------------------------------
# Copy row 102's A-cell first so the new date cell inherits the existing
# date style (s="1") instead of Excel creating a brand new style entry.
$ws.Range("A102").Copy($ws.Range("A103"))
$ws.Cells.Item(103, 1).Value = 45490.2916666667
$ws.Cells.Item(103, 2).Value = 0
$ws.Cells.Item(103, 3).Value = 3.29999995231628
$ws.Cells.Item(103, 4).Value = 3.29999995231628
$ws.Cells.Item(103, 5).Value = 3.29999995231628
$ws.Cells.Item(103, 6).Value = 3.29999995231628

# adj_close is stored as text in this sheet (matches the "3.29999995231628"
# shared string that already exists for the other rows with that close
# price). Force text typing via the "@" format, then strip the formatting
# back off so the cell keeps the default style, same as its neighbours.
$ws.Cells.Item(103, 7).NumberFormat = "@"
$ws.Cells.Item(103, 7).Value = "3.29999995231628"
$ws.Cells.Item(103, 7).ClearFormats()

$ws.Cells.Item(103, 8).Value = "ESPE.MI"

# --- Row 104 -----------------------------------------------------------
$ws.Range("A102").Copy($ws.Range("A104"))
$ws.Cells.Item(104, 1).Value = 45491.5507060185
$ws.Cells.Item(104, 2).Value = 6000
$ws.Cells.Item(104, 3).Value = 3.24000000953674
$ws.Cells.Item(104, 4).Value = 3.15000009536743
$ws.Cells.Item(104, 5).Value = 3.1800000667572
$ws.Cells.Item(104, 6).Value = 3.21000003814697

$ws.Cells.Item(104, 7).NumberFormat = "@"
$ws.Cells.Item(104, 7).Value = "3.21000003814697"
$ws.Cells.Item(104, 7).ClearFormats()

$ws.Cells.Item(104, 8).Value = "ESPE.MI"
